$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $rng1 = $ws.Range("B$row1`:AD$row1")
    $rng2 = $ws.Range("B$row2`:AD$row2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Swap-Rows 186 187
Swap-Rows 243 244
Swap-Rows 260 261
Swap-Rows 296 297
